$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") holds values such as "System, dnasr281@gmail.com".
# Swap the order of the two names to "dnasr281@gmail.com, System" wherever
# that exact value occurs, leaving every other cell (including other G-column
# values like a lone "System" or a lone email) untouched.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
